# Scheduled market-data refresh: update currentAveragePrice* / Leve* price &
# profit columns (H:N) across the Leve-profit tracker sheets with freshly
# fetched values. Plain data values only - no formulas involved.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ALC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H15").Value = 1608.3334
$ws.Range("I15").Value = 1608.3334
$ws.Range("K15").Value = 4825.0002
$ws.Range("M15").Value = -4656.0002

$ws.Range("H33").Value = 690703.5600000001
$ws.Range("I33").Value = 958421.5600000001
$ws.Range("K33").Value = 958421.5600000001
$ws.Range("M33").Value = -958192.5600000001

$ws.Range("H53").Value = 155
$ws.Range("I53").Value = 151.38461
$ws.Range("J53").Value = 157.2381
$ws.Range("K53").Value = 151.38461
$ws.Range("L53").Value = 157.2381
$ws.Range("M53").Value = 485.61539
$ws.Range("N53").Value = -1431.2381

$ws.Range("H62").Value = 25413.818
$ws.Range("I62").Value = 20000.285
$ws.Range("J62").Value = 34887.5
$ws.Range("K62").Value = 20000.285
$ws.Range("L62").Value = 34887.5
$ws.Range("M62").Value = -19376.285
$ws.Range("N62").Value = -36135.5

$ws.Range("H65").Value = 25413.818
$ws.Range("I65").Value = 20000.285
$ws.Range("J65").Value = 34887.5
$ws.Range("K65").Value = 100001.425
$ws.Range("L65").Value = 174437.5
$ws.Range("M65").Value = -96881.425
$ws.Range("N65").Value = -180677.5

$ws.Range("H68").Value = 25221
$ws.Range("J68").Value = 25294.666
$ws.Range("L68").Value = 25294.666
$ws.Range("N68").Value = -26792.666

$ws.Range("H70").Value = 12105.889
$ws.Range("I70").Value = 17617
$ws.Range("J70").Value = 1083.6666
$ws.Range("K70").Value = 52851
$ws.Range("L70").Value = 3250.9998
$ws.Range("M70").Value = -52581
$ws.Range("N70").Value = -3790.9998

$ws.Range("H71").Value = 25221
$ws.Range("J71").Value = 25294.666
$ws.Range("L71").Value = 75883.99800000001
$ws.Range("N71").Value = -83371.99800000001

$ws.Range("H73").Value = 12105.889
$ws.Range("I73").Value = 17617
$ws.Range("J73").Value = 1083.6666
$ws.Range("K73").Value = 52851
$ws.Range("L73").Value = 3250.9998
$ws.Range("M73").Value = -51915
$ws.Range("N73").Value = -5122.9998

$ws.Range("H86").Value = 3994.5
$ws.Range("I86").Value = 3994.5
$ws.Range("K86").Value = 3994.5
$ws.Range("M86").Value = -2871.5

$ws.Range("H89").Value = 3994.5
$ws.Range("I89").Value = 3994.5
$ws.Range("K89").Value = 19972.5
$ws.Range("M89").Value = -14356.5

$ws.Range("H106").Value = 2412.111
$ws.Range("I106").Value = 1643.2
$ws.Range("K106").Value = 1643.2
$ws.Range("M106").Value = -1012.2

$ws.Range("H112").Value = 1967.5
$ws.Range("J112").Value = 1719.091
$ws.Range("L112").Value = 5157.272999999999
$ws.Range("N112").Value = -7373.272999999999

$ws.Range("H132").Value = 3444.8333
$ws.Range("I132").Value = 1487.44
$ws.Range("K132").Value = 4462.32
$ws.Range("M132").Value = -1932.32

$ws.Range("H137").Value = 3276.0938
$ws.Range("I137").Value = 2544.9375
$ws.Range("J137").Value = 4007.25
$ws.Range("K137").Value = 7634.8125
$ws.Range("L137").Value = 12021.75
$ws.Range("M137").Value = -5084.8125
$ws.Range("N137").Value = -17121.75

$ws.Range("H138").Value = 2565.0193
$ws.Range("I138").Value = 1732.3
$ws.Range("J138").Value = 3085.4688
$ws.Range("K138").Value = 5196.9
$ws.Range("L138").Value = 9256.4064
$ws.Range("M138").Value = -56.89999999999964
$ws.Range("N138").Value = -19536.4064

# ---------------------------------------------------------------------------
# ARM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("I32").Value = 816854.0600000001
$ws.Range("K32").Value = 816854.0600000001
$ws.Range("M32").Value = -816567.0600000001

$ws.Range("H34").Value = 16859
$ws.Range("I34").Value = 16859
$ws.Range("K34").Value = 16859
$ws.Range("M34").Value = -16588

$ws.Range("H43").Value = 6168.75
$ws.Range("J43").Value = 6168.75
$ws.Range("L43").Value = 6168.75
$ws.Range("N43").Value = -6794.75

$ws.Range("H45").Value = 35191.777
$ws.Range("I45").Value = 27089
$ws.Range("K45").Value = 27089
$ws.Range("M45").Value = -26712

# ---------------------------------------------------------------------------
# BSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H134").Value = 3795.7273
$ws.Range("I134").Value = 3232.9312
$ws.Range("K134").Value = 9698.793600000001
$ws.Range("M134").Value = -7163.793600000001

# ---------------------------------------------------------------------------
# CRP
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 3963.8333
$ws.Range("I31").Value = 2249.75
$ws.Range("J31").Value = 7392
$ws.Range("K31").Value = 2249.75
$ws.Range("L31").Value = 7392
$ws.Range("M31").Value = -1954.75
$ws.Range("N31").Value = -7982

$ws.Range("H34").Value = 3963.8333
$ws.Range("I34").Value = 2249.75
$ws.Range("J34").Value = 7392
$ws.Range("K34").Value = 2249.75
$ws.Range("L34").Value = 7392
$ws.Range("M34").Value = -2047.75
$ws.Range("N34").Value = -7796

$ws.Range("H94").Value = 2545.923
$ws.Range("I94").Value = 2797.8
$ws.Range("K94").Value = 2797.8
$ws.Range("M94").Value = -2346.8

$ws.Range("H132").Value = 5314.857
$ws.Range("I132").Value = 4761
$ws.Range("K132").Value = 14283
$ws.Range("M132").Value = -11753

$ws.Range("H141").Value = 130333.336
$ws.Range("J141").Value = 130333.336
$ws.Range("L141").Value = 130333.336
$ws.Range("N141").Value = -140693.336

# ---------------------------------------------------------------------------
# CUL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H50").Value = 889191.6
$ws.Range("I50").Value = 100115.8
$ws.Range("J50").Value = 1875536.4
$ws.Range("K50").Value = 300347.4
$ws.Range("L50").Value = 5626609.199999999
$ws.Range("M50").Value = -299866.4
$ws.Range("N50").Value = -5627571.199999999

$ws.Range("H53").Value = 889191.6
$ws.Range("I53").Value = 100115.8
$ws.Range("J53").Value = 1875536.4
$ws.Range("K53").Value = 300347.4
$ws.Range("L53").Value = 5626609.199999999
$ws.Range("M53").Value = -299866.4
$ws.Range("N53").Value = -5627571.199999999

# ---------------------------------------------------------------------------
# GSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H70").Value = 37646850
$ws.Range("I70").Value = 5490.143
$ws.Range("K70").Value = 5490.143
$ws.Range("M70").Value = -5220.143

$ws.Range("H73").Value = 37646850
$ws.Range("I73").Value = 5490.143
$ws.Range("K73").Value = 5490.143
$ws.Range("M73").Value = -4554.143

# Row 102: M keeps the (new) value that used to live in N; N's cell goes away.
$ws.Range("H102").Value = 2958.1667
$ws.Range("I102").Value = 2958.1667
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2958.1667
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("M102").Value = -1336.1667

# Row 122: N's cell goes away entirely (M stays at its existing value).
$ws.Range("H122").Value = 2306.1538
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 311571
$ws.Range("J126").Value = 311571
$ws.Range("L126").Value = 934713
$ws.Range("N126").Value = -939653

$ws.Range("H132").Value = 4255.1
$ws.Range("I132").Value = 4362
$ws.Range("J132").Value = 3827.5
$ws.Range("K132").Value = 13086
$ws.Range("L132").Value = 11482.5
$ws.Range("M132").Value = -10556
$ws.Range("N132").Value = -16542.5

# ---------------------------------------------------------------------------
# LTW
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

# Row 40: M's cell goes away entirely; N keeps the (new) value that used to
# live in M.
$ws.Range("H40").Value = 9924.75
$ws.Range("I40").Value = 9924.75
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 9924.75
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -9788.75

$ws.Range("H82").Value = 822.0909
$ws.Range("I82").Value = 673.1111
$ws.Range("K82").Value = 673.1111
$ws.Range("M82").Value = -312.1111

$ws.Range("H85").Value = 822.0909
$ws.Range("I85").Value = 673.1111
$ws.Range("K85").Value = 673.1111
$ws.Range("M85").Value = 574.8889

# ---------------------------------------------------------------------------
# WVR
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# Row 40: L's cell goes away entirely; N keeps the (new) value.
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("N40").Value = 0

# Row 62: L's cell goes away entirely; M keeps the (new) value that used to
# live in L.
$ws.Range("H62").Value = 170453.81
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 170453.81
$ws.Range("K62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("M62").Value = 170453.81
$ws.Range("N62").Value = -171701.81

# Row 65: same pattern as row 62.
$ws.Range("H65").Value = 170453.81
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 170453.81
$ws.Range("K65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("M65").Value = 852269.05
$ws.Range("N65").Value = -858509.05

$ws.Range("H81").Value = 11761.75
$ws.Range("I81").Value = 1365
$ws.Range("J81").Value = 17999.8
$ws.Range("K81").Value = 2730
$ws.Range("L81").Value = 35999.6
$ws.Range("M81").Value = -1669
$ws.Range("N81").Value = -38121.6

$ws.Range("H84").Value = 11761.75
$ws.Range("I84").Value = 1365
$ws.Range("J84").Value = 17999.8
$ws.Range("K84").Value = 13650
$ws.Range("L84").Value = 179998
$ws.Range("M84").Value = -8346
$ws.Range("N84").Value = -190606

$ws.Range("H132").Value = 7907
$ws.Range("I132").Value = 8335.044
$ws.Range("K132").Value = 25005.132
$ws.Range("M132").Value = -22475.132

$ws.Range("H136").Value = 7273.5186
$ws.Range("I136").Value = 4860.769
$ws.Range("J136").Value = 70005
$ws.Range("K136").Value = 14582.307
$ws.Range("L136").Value = 210015
$ws.Range("M136").Value = -12032.307
$ws.Range("N136").Value = -215115
